$d = $word.ActiveDocument

# 1. Change the "Obra" definition text
$d.Content.Find.Execute(
    ": Bien inmueble que será construido o reformado.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": Información referente a un proyecto.", 2
)

# 2. Change heading "Unidad Habitacional (UH)" -> "Inmueble"
$d.Content.Find.Execute(
    "Unidad Habitacional (UH)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Inmueble", 2
)

# 3. Change the text right after the heading, up to the first "Obra"
$d.Content.Find.Execute(
    ": Cada una de las subdivisiones de una ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": El producto terminado que será construido o reformado. Ej: en la ",
    2
)

# 4. Change the text between the two bold "Obra" runs
$d.Content.Find.Execute(
    ". Ej: en la Obra de un edificio, cada UH es un departamento. En la ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " de un edificio, cada inmueble es un departamento. En la ",
    2
)

# 5. Change the final trailing text
$d.Content.Find.Execute(
    " de un barrio, cada UH es una vivienda.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " de un barrio, es una casa.",
    2
)
